$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add column E formulas "=D{row}-18" on every row that already has a
# weekly SUM() formula in column D.
$weekRows = @(5, 9, 13, 17, 22, 25, 29, 34, 37, 42, 46, 52)
foreach ($r in $weekRows) {
    $ws.Range("E$r").Formula = "=D$r-18"
}

# Add the two new data rows (55 and 56) that were previously blank.
$ws.Range("A55").Value = 41472
$ws.Range("B55").Value = "Handbuch, Aspekt 3 Commands / Handler, Menus, Key Bindings"
$ws.Range("C55").Value = 1

$ws.Range("A56").Value = 41473
$ws.Range("B56").Value = "Handbuch, Aspekt 3 Commands / Handler, Menus, Key Bindings"
$ws.Range("C56").Value = 8

# Update the view so it scrolls to roughly where the new data is and the
# active selection is on C57 (the next empty data cell).
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("C57").Select()
